$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demo")

# Add a new row of data (row 6) below the existing data
$ws.Range("A6").Value = "Modify profile"
$ws.Range("B6").Value = "mmm"
$ws.Range("C6").Value = "nnn"
$ws.Range("D6").Value = "ooo"

# Update selection to match the diff (active cell B7)
$ws.Range("B7").Select()
